$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: "Pull Request 3 :" section header -- same green fill as the
# "Pull Request 2 :" banner (A10) combined with the plain black RGB font
# already used elsewhere in the sheet (A1/A2).
$ws.Range("A14").Value = "Pull Request 3 :"
$ws.Range("A14").Font.Color = $ws.Range("A2").Font.Color
$ws.Range("A14").Interior.Color = $ws.Range("A10").Interior.Color

# Row 15: git branch command + its Google Drive link
$ws.Range("A15").Value = "git branch command"
$ws.Range("C15").Value = "https://drive.google.com/file/d/1EicC-7Y0EWyOC_GZoFwkSCuFiFskn6Hi/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("C15"), "https://drive.google.com/file/d/1EicC-7Y0EWyOC_GZoFwkSCuFiFskn6Hi/view?usp=sharing")
$ws.Range("C15").Style = "Hyperlink"

# Row 16: git log command + its Google Drive link
$ws.Range("A16").Value = "git log command"
$ws.Range("C16").Value = "https://drive.google.com/file/d/1kjxST2z_KyeUqAHRPN4_hn6SFRp-Ye8C/view?usp=sharing"
$ws.Hyperlinks.Add($ws.Range("C16"), "https://drive.google.com/file/d/1kjxST2z_KyeUqAHRPN4_hn6SFRp-Ye8C/view?usp=sharing")
$ws.Range("C16").Style = "Hyperlink"

# New rows render a bit taller (15pt) than the sheet's default row height.
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 15
$ws.Rows.Item(16).RowHeight = 15

# Leave the cursor where the author ended up after typing the new rows.
$ws.Range("B16").Select()
